# feat: add 2022-Q1 data
#
# The workbook tracks per-quarter fund-holdings snapshots (one sheet per
# quarter) plus a rolling "总计" (totals) summary sheet.
#
# This change:
#   1. Turns the current "总计" sheet into the new "2022-Q1" holdings sheet
#      (it keeps its sheetId/position, it is simply renamed + refilled with
#      the 2022-Q1 per-fund holdings data).
#   2. Appends a brand new "总计" sheet at the end of the workbook with the
#      summary table, now including a new leading row for 2022-Q1.

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# Step 1: repurpose the existing "总计" sheet into the "2022-Q1" sheet
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

# Remove the old 5th data row (row 6 held the "2020-Q4" summary record);
# the new sheet only needs 4 fund rows (rows 2-5).
$q1.Rows(6).Delete()

# -- helper: write a value as literal text (no auto number conversion) --
function Set-TextValue($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# Headers (row 1). B1:D1 already carry the bold/bordered header style from
# the previous "总计" content; just overwrite the captions and extend the
# same style across the new E1:H1 columns.
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"

$q1.Range("D1").Copy()
$q1.Range("E1:H1").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Data rows 2-5 (column A index values 0,1,2,3 are already correct and are
# left untouched).
Set-TextValue $q1.Range("B2") "010783"
Set-TextValue $q1.Range("C2") "德邦沪港深龙头混合A"
Set-TextValue $q1.Range("D2") "0.93"
Set-TextValue $q1.Range("E2") "81.58"
Set-TextValue $q1.Range("F2") "5.38"
Set-TextValue $q1.Range("G2") "0.0500"
$q1.Range("H2").Value = 5

Set-TextValue $q1.Range("B3") "513160"
Set-TextValue $q1.Range("C3") "银华恒生港股通中国科技ETF"
Set-TextValue $q1.Range("D3") "0.62"
Set-TextValue $q1.Range("E3") "92.07"
Set-TextValue $q1.Range("F3") "4.55"
Set-TextValue $q1.Range("G3") "0.0282"
$q1.Range("H3").Value = 9

Set-TextValue $q1.Range("B4") "010784"
Set-TextValue $q1.Range("C4") "德邦沪港深龙头混合C"
Set-TextValue $q1.Range("D4") "0.27"
Set-TextValue $q1.Range("E4") "81.58"
Set-TextValue $q1.Range("F4") "5.38"
Set-TextValue $q1.Range("G4") "0.0145"
$q1.Range("H4").Value = 5

Set-TextValue $q1.Range("B5") "160922"
Set-TextValue $q1.Range("C5") "大成恒生综合中小型股指数(QDII-LOF)A"
Set-TextValue $q1.Range("D5") "0.10"
Set-TextValue $q1.Range("E5") "92.44"
Set-TextValue $q1.Range("F5") "1.04"
Set-TextValue $q1.Range("G5") "0.0010"
$q1.Range("H5").Value = 9

# ---------------------------------------------------------------------
# Step 2: append a fresh "总计" sheet at the end with the updated summary
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$total = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$total.Name = "总计"

# Match the look & feel (page margins, outline defaults, ...) of the other
# data sheets.
$total.PageSetup.LeftMargin = 54
$total.PageSetup.RightMargin = 54
$total.PageSetup.TopMargin = 72
$total.PageSetup.BottomMargin = 72
$total.PageSetup.HeaderMargin = 36
$total.PageSetup.FooterMargin = 36
$total.Outline.SummaryRow = 1
$total.Outline.SummaryColumn = 1

# Copy the header style (bold + border) from one of the existing sheets.
$q1.Range("B1").Copy()
$total.Range("B1:D1").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

# Copy the column-A index style too.
$q1.Range("A2").Copy()
$total.Range("A2:A7").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

$rows = @(
    @(0, "2022-Q1", 4, 0.09),
    @(1, "2021-Q4", 2, 0.01),
    @(2, "2021-Q3", 2, 0.01),
    @(3, "2021-Q2", 2, 0.01),
    @(4, "2021-Q1", 2, 0.01),
    @(5, "2020-Q4", 2, 0.01)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $total.Range("A$r").Value = $row[0]
    $total.Range("B$r").Value = $row[1]
    $total.Range("C$r").Value = $row[2]
    $total.Range("D$r").Value = $row[3]
}

# Restore the originally active sheet / selection.
$wb.Worksheets.Item(1).Activate()
